# Commit: "add image and publish documents"
#
# The only substantive content edit is on slide 1: the (until now
# empty) Subtitle placeholder gets a run of text "eyreryrty", authored
# with the deck's existing "en-IN" language tag (matching the sibling
# Title placeholder's run, rather than a generic "en-US" default).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$subtitle = $s.Shapes.Item(2)
$tr = $subtitle.TextFrame.TextRange
$tr.Text = "eyreryrty"
$tr.LanguageID = "en-IN"

# Best-effort/no-op elsewhere: the author's save also picked up an
# (empty) PowerPoint 2013+ slide-guide list on the presentation itself
# (p:extLst/p15:sldGuideLst). Touch the Guides collection so this is
# picked up if/when the host models it; harmless otherwise.
try {
    $null = $p.Guides
} catch {
}
